$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 266-290 with the corrected weekly price data
$ws.Cells.Item(266, 4).Value = 45180
$ws.Cells.Item(266, 12).Value = "Especial"
$ws.Cells.Item(266, 13).Value = 200
$ws.Cells.Item(266, 14).Value = 21000
$ws.Cells.Item(266, 15).Value = 22000
$ws.Cells.Item(266, 16).Value = 21500
$ws.Cells.Item(266, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(266, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(266, 19).Value = 2150
$ws.Cells.Item(266, 20).Value = 10

$ws.Cells.Item(267, 4).Value = 45180
$ws.Cells.Item(267, 12).Value = "Primera"
$ws.Cells.Item(267, 13).Value = 240
$ws.Cells.Item(267, 14).Value = 19000
$ws.Cells.Item(267, 15).Value = 20000
$ws.Cells.Item(267, 16).Value = 19500
$ws.Cells.Item(267, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(267, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(267, 19).Value = 1950
$ws.Cells.Item(267, 20).Value = 10

$ws.Cells.Item(268, 4).Value = 45180
$ws.Cells.Item(268, 12).Value = "Segunda"
$ws.Cells.Item(268, 13).Value = 240
$ws.Cells.Item(268, 14).Value = 16000
$ws.Cells.Item(268, 15).Value = 17000
$ws.Cells.Item(268, 16).Value = 16500
$ws.Cells.Item(268, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(268, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(268, 19).Value = 1650
$ws.Cells.Item(268, 20).Value = 10

$ws.Cells.Item(269, 4).Value = 45166
$ws.Cells.Item(269, 12).Value = "Especial"
$ws.Cells.Item(269, 13).Value = 160
$ws.Cells.Item(269, 14).Value = 24000
$ws.Cells.Item(269, 15).Value = 25000
$ws.Cells.Item(269, 16).Value = 24500
$ws.Cells.Item(269, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(269, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(269, 19).Value = 2450
$ws.Cells.Item(269, 20).Value = 10

$ws.Cells.Item(270, 4).Value = 45166
$ws.Cells.Item(270, 12).Value = "Primera"
$ws.Cells.Item(270, 13).Value = 240
$ws.Cells.Item(270, 14).Value = 22000
$ws.Cells.Item(270, 15).Value = 23000
$ws.Cells.Item(270, 16).Value = 22500
$ws.Cells.Item(270, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(270, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(270, 19).Value = 2250
$ws.Cells.Item(270, 20).Value = 10

$ws.Cells.Item(271, 4).Value = 45166
$ws.Cells.Item(271, 12).Value = "Segunda"
$ws.Cells.Item(271, 13).Value = 240
$ws.Cells.Item(271, 14).Value = 18000
$ws.Cells.Item(271, 15).Value = 19000
$ws.Cells.Item(271, 16).Value = 18500
$ws.Cells.Item(271, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(271, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(271, 19).Value = 1850
$ws.Cells.Item(271, 20).Value = 10

$ws.Cells.Item(272, 4).Value = 44459
$ws.Cells.Item(272, 12).Value = "Especial"
$ws.Cells.Item(272, 13).Value = 200
$ws.Cells.Item(272, 14).Value = 2600
$ws.Cells.Item(272, 15).Value = 2700
$ws.Cells.Item(272, 16).Value = 2650
$ws.Cells.Item(272, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(272, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(272, 19).Value = 2650
$ws.Cells.Item(272, 20).Value = 1

$ws.Cells.Item(273, 4).Value = 44459
$ws.Cells.Item(273, 12).Value = "Primera"
$ws.Cells.Item(273, 13).Value = 300
$ws.Cells.Item(273, 14).Value = 2200
$ws.Cells.Item(273, 15).Value = 2300
$ws.Cells.Item(273, 16).Value = 2250
$ws.Cells.Item(273, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(273, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(273, 19).Value = 2250
$ws.Cells.Item(273, 20).Value = 1

$ws.Cells.Item(274, 4).Value = 44459
$ws.Cells.Item(274, 12).Value = "Segunda"
$ws.Cells.Item(274, 13).Value = 240
$ws.Cells.Item(274, 14).Value = 1900
$ws.Cells.Item(274, 15).Value = 2000
$ws.Cells.Item(274, 16).Value = 1950
$ws.Cells.Item(274, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(274, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(274, 19).Value = 1950
$ws.Cells.Item(274, 20).Value = 1

$ws.Cells.Item(275, 4).Value = 44880
$ws.Cells.Item(275, 12).Value = "Especial"
$ws.Cells.Item(275, 13).Value = 360
$ws.Cells.Item(275, 14).Value = 18000
$ws.Cells.Item(275, 15).Value = 19000
$ws.Cells.Item(275, 16).Value = 18500
$ws.Cells.Item(275, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(275, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(275, 19).Value = 1850
$ws.Cells.Item(275, 20).Value = 10

$ws.Cells.Item(276, 4).Value = 44880
$ws.Cells.Item(276, 12).Value = "Primera"
$ws.Cells.Item(276, 13).Value = 300
$ws.Cells.Item(276, 14).Value = 15000
$ws.Cells.Item(276, 15).Value = 16000
$ws.Cells.Item(276, 16).Value = 15500
$ws.Cells.Item(276, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(276, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(276, 19).Value = 1550
$ws.Cells.Item(276, 20).Value = 10

$ws.Cells.Item(277, 4).Value = 44880
$ws.Cells.Item(277, 12).Value = "Segunda"
$ws.Cells.Item(277, 13).Value = 200
$ws.Cells.Item(277, 14).Value = 12000
$ws.Cells.Item(277, 15).Value = 13000
$ws.Cells.Item(277, 16).Value = 12500
$ws.Cells.Item(277, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(277, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(277, 19).Value = 1250
$ws.Cells.Item(277, 20).Value = 10

$ws.Cells.Item(278, 4).Value = 44516
$ws.Cells.Item(278, 12).Value = "Especial"
$ws.Cells.Item(278, 13).Value = 360
$ws.Cells.Item(278, 14).Value = 1800
$ws.Cells.Item(278, 15).Value = 1900
$ws.Cells.Item(278, 16).Value = 1850
$ws.Cells.Item(278, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(278, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(278, 19).Value = 1850
$ws.Cells.Item(278, 20).Value = 1

$ws.Cells.Item(279, 4).Value = 44516
$ws.Cells.Item(279, 12).Value = "Primera"
$ws.Cells.Item(279, 13).Value = 320
$ws.Cells.Item(279, 14).Value = 1500
$ws.Cells.Item(279, 15).Value = 1600
$ws.Cells.Item(279, 16).Value = 1550
$ws.Cells.Item(279, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(279, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(279, 19).Value = 1550
$ws.Cells.Item(279, 20).Value = 1

$ws.Cells.Item(280, 4).Value = 44516
$ws.Cells.Item(280, 12).Value = "Segunda"
$ws.Cells.Item(280, 13).Value = 240
$ws.Cells.Item(280, 14).Value = 1200
$ws.Cells.Item(280, 15).Value = 1300
$ws.Cells.Item(280, 16).Value = 1250
$ws.Cells.Item(280, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(280, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(280, 19).Value = 1250
$ws.Cells.Item(280, 20).Value = 1

$ws.Cells.Item(281, 4).Value = 44848
$ws.Cells.Item(281, 12).Value = "Especial"
$ws.Cells.Item(281, 13).Value = 240
$ws.Cells.Item(281, 14).Value = 21000
$ws.Cells.Item(281, 15).Value = 22000
$ws.Cells.Item(281, 16).Value = 21500
$ws.Cells.Item(281, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(281, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(281, 19).Value = 2150
$ws.Cells.Item(281, 20).Value = 10

$ws.Cells.Item(282, 4).Value = 44848
$ws.Cells.Item(282, 12).Value = "Primera"
$ws.Cells.Item(282, 13).Value = 400
$ws.Cells.Item(282, 14).Value = 18000
$ws.Cells.Item(282, 15).Value = 19000
$ws.Cells.Item(282, 16).Value = 18500
$ws.Cells.Item(282, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(282, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(282, 19).Value = 1850
$ws.Cells.Item(282, 20).Value = 10

$ws.Cells.Item(283, 4).Value = 44848
$ws.Cells.Item(283, 12).Value = "Segunda"
$ws.Cells.Item(283, 13).Value = 400
$ws.Cells.Item(283, 14).Value = 15000
$ws.Cells.Item(283, 15).Value = 16000
$ws.Cells.Item(283, 16).Value = 15500
$ws.Cells.Item(283, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(283, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(283, 19).Value = 1550
$ws.Cells.Item(283, 20).Value = 10

$ws.Cells.Item(284, 4).Value = 44530
$ws.Cells.Item(284, 12).Value = "Especial"
$ws.Cells.Item(284, 13).Value = 440
$ws.Cells.Item(284, 14).Value = 1700
$ws.Cells.Item(284, 15).Value = 1800
$ws.Cells.Item(284, 16).Value = 1750
$ws.Cells.Item(284, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(284, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(284, 19).Value = 1750
$ws.Cells.Item(284, 20).Value = 1

$ws.Cells.Item(285, 4).Value = 44530
$ws.Cells.Item(285, 12).Value = "Primera"
$ws.Cells.Item(285, 13).Value = 400
$ws.Cells.Item(285, 14).Value = 1500
$ws.Cells.Item(285, 15).Value = 1600
$ws.Cells.Item(285, 16).Value = 1550
$ws.Cells.Item(285, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(285, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(285, 19).Value = 1550
$ws.Cells.Item(285, 20).Value = 1

$ws.Cells.Item(286, 4).Value = 44530
$ws.Cells.Item(286, 12).Value = "Segunda"
$ws.Cells.Item(286, 13).Value = 280
$ws.Cells.Item(286, 14).Value = 1200
$ws.Cells.Item(286, 15).Value = 1300
$ws.Cells.Item(286, 16).Value = 1250
$ws.Cells.Item(286, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(286, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(286, 19).Value = 1250
$ws.Cells.Item(286, 20).Value = 1

$ws.Cells.Item(287, 4).Value = 44813
$ws.Cells.Item(287, 12).Value = "Especial"
$ws.Cells.Item(287, 13).Value = 300
$ws.Cells.Item(287, 14).Value = 23000
$ws.Cells.Item(287, 15).Value = 24000
$ws.Cells.Item(287, 16).Value = 23500
$ws.Cells.Item(287, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(287, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(287, 19).Value = 2350
$ws.Cells.Item(287, 20).Value = 10

$ws.Cells.Item(288, 4).Value = 44813
$ws.Cells.Item(288, 12).Value = "Primera"
$ws.Cells.Item(288, 13).Value = 400
$ws.Cells.Item(288, 14).Value = 19000
$ws.Cells.Item(288, 15).Value = 20000
$ws.Cells.Item(288, 16).Value = 19500
$ws.Cells.Item(288, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(288, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(288, 19).Value = 1950
$ws.Cells.Item(288, 20).Value = 10

$ws.Cells.Item(289, 4).Value = 44813
$ws.Cells.Item(289, 12).Value = "Segunda"
$ws.Cells.Item(289, 13).Value = 400
$ws.Cells.Item(289, 14).Value = 16000
$ws.Cells.Item(289, 15).Value = 17000
$ws.Cells.Item(289, 16).Value = 16500
$ws.Cells.Item(289, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(289, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(289, 19).Value = 1650
$ws.Cells.Item(289, 20).Value = 10

$ws.Cells.Item(290, 4).Value = 44813
$ws.Cells.Item(290, 12).Value = "Tercera"
$ws.Cells.Item(290, 13).Value = 300
$ws.Cells.Item(290, 14).Value = 13000
$ws.Cells.Item(290, 15).Value = 14000
$ws.Cells.Item(290, 16).Value = 13500
$ws.Cells.Item(290, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(290, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(290, 19).Value = 1350
$ws.Cells.Item(290, 20).Value = 10

# Add new rows 291-293 (previous weekly entries pushed down)
$ws.Cells.Item(291, 1).Value = 8
$ws.Cells.Item(291, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(291, 3).Value = "Coquimbo"
$ws.Cells.Item(291, 4).Value = 44490
$ws.Cells.Item(291, 5).Value = 4
$ws.Cells.Item(291, 6).Value = "Fruta"
$ws.Cells.Item(291, 7).Value = 100107
$ws.Cells.Item(291, 8).Value = "Otros"
$ws.Cells.Item(291, 9).Value = 100107002
$ws.Cells.Item(291, 10).Value = "Chirimoya"
$ws.Cells.Item(291, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(291, 12).Value = "Especial"
$ws.Cells.Item(291, 13).Value = 240
$ws.Cells.Item(291, 14).Value = 2200
$ws.Cells.Item(291, 15).Value = 2300
$ws.Cells.Item(291, 16).Value = 2250
$ws.Cells.Item(291, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(291, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(291, 19).Value = 2250
$ws.Cells.Item(291, 20).Value = 1

$ws.Cells.Item(292, 1).Value = 8
$ws.Cells.Item(292, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(292, 3).Value = "Coquimbo"
$ws.Cells.Item(292, 4).Value = 44490
$ws.Cells.Item(292, 5).Value = 4
$ws.Cells.Item(292, 6).Value = "Fruta"
$ws.Cells.Item(292, 7).Value = 100107
$ws.Cells.Item(292, 8).Value = "Otros"
$ws.Cells.Item(292, 9).Value = 100107002
$ws.Cells.Item(292, 10).Value = "Chirimoya"
$ws.Cells.Item(292, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(292, 12).Value = "Primera"
$ws.Cells.Item(292, 13).Value = 500
$ws.Cells.Item(292, 14).Value = 1900
$ws.Cells.Item(292, 15).Value = 2000
$ws.Cells.Item(292, 16).Value = 1950
$ws.Cells.Item(292, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(292, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(292, 19).Value = 1950
$ws.Cells.Item(292, 20).Value = 1

$ws.Cells.Item(293, 1).Value = 8
$ws.Cells.Item(293, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(293, 3).Value = "Coquimbo"
$ws.Cells.Item(293, 4).Value = 44490
$ws.Cells.Item(293, 5).Value = 4
$ws.Cells.Item(293, 6).Value = "Fruta"
$ws.Cells.Item(293, 7).Value = 100107
$ws.Cells.Item(293, 8).Value = "Otros"
$ws.Cells.Item(293, 9).Value = 100107002
$ws.Cells.Item(293, 10).Value = "Chirimoya"
$ws.Cells.Item(293, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(293, 12).Value = "Segunda"
$ws.Cells.Item(293, 13).Value = 360
$ws.Cells.Item(293, 14).Value = 1400
$ws.Cells.Item(293, 15).Value = 1500
$ws.Cells.Item(293, 16).Value = 1450
$ws.Cells.Item(293, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(293, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(293, 19).Value = 1450
$ws.Cells.Item(293, 20).Value = 1

# Match the date-time number format used by the other date cells in column D
$ws.Range("D291:D293").NumberFormat = "YYYY-MM-DD HH:MM:SS"
